$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 206
$ws.Range("B2").Value = 14
$ws.Range("C2").Value = "許*綸"
$ws.Range("D2").Value = "2024-03-01 12:59:50"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "IN"

$ws.Range("A3").Value = 206
$ws.Range("B3").Value = 14
$ws.Range("C3").Value = "許*綸"
$ws.Range("D3").Value = "2024-03-01 12:59:39"
$ws.Range("E3").Value = "2024-03-01 12:59:47"
$ws.Range("F3").Value = "OUT"
